# Add carjacking data for 2022-06-05 (advance the "through" date by one day)
# to the "carjacking-by-neighborhood-by-month" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the column-header label to reflect the new
# "through" date.
$ws.Name = "Through 2022-06-05"
$ws.Range("B1").Value = "June 2022 (through June 05)"

# Englewood (row 2): June 2021 count 2 -> 3
$ws.Range("H2").Value = 3

# North Lawndale (row 4): June 2022 count 1 -> 2
$ws.Range("B4").Value = 2

# South Shore (row 5): new incident in June 2016
$ws.Range("AL5").Value = 1

# Humboldt Park (row 6): new incident in June 2022
$ws.Range("B6").Value = 1

# Logan Square (row 8): new incident in June 2018
$ws.Range("Z8").Value = 1

# Grand Crossing (row 9): new incident in June 2022
$ws.Range("B9").Value = 1

# Garfield Park (row 10): June 2020 count 1 -> 2; new incident June 2016
$ws.Range("N10").Value = 2
$ws.Range("AL10").Value = 1

# Washington Park (row 11): new incident in June 2021
$ws.Range("H11").Value = 1

# Lincoln Park (row 21): new incident in June 2021
$ws.Range("H21").Value = 1

# West Ridge (row 29): new incident in June 2021
$ws.Range("H29").Value = 1

# Avondale (row 46): new incident in June 2018
$ws.Range("Z46").Value = 1

# Irving Park (row 64): June 2020 count 1 -> 2
$ws.Range("N64").Value = 2

# Jackson Park (row 65): new incident in June 2022
$ws.Range("B65").Value = 1

# River North (row 84): new incident in June 2021
$ws.Range("H84").Value = 1

# West Pullman (row 93): new incident in June 2019
$ws.Range("T93").Value = 1
